$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared string rich text runs) ---
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Data table updates (rows 14-30) ---
$ws.Range("D14").Value = 1
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = -54.545454545454
$ws.Range("J14").Value = 65
$ws.Range("K14").Value = -13.846153846153
$ws.Range("L14").Value = -30
$ws.Range("M14").Value = -51.724137931034
$ws.Range("N14").Value = -85.858585858585
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = -22.222222222222
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 26
$ws.Range("H15").Value = -23.076923076923
$ws.Range("I15").Value = 185
$ws.Range("J15").Value = 214
$ws.Range("K15").Value = -13.551401869158
$ws.Range("L15").Value = 2.209944751381
$ws.Range("M15").Value = 4.519774011299
$ws.Range("N15").Value = -62.321792260692
$ws.Range("C16").Value = 64
$ws.Range("D16").Value = 44
$ws.Range("E16").Value = 45.454545454545
$ws.Range("F16").Value = 210
$ws.Range("G16").Value = 182
$ws.Range("H16").Value = 15.384615384615
$ws.Range("I16").Value = 2004
$ws.Range("J16").Value = 2101
$ws.Range("K16").Value = -4.616849119466
$ws.Range("L16").Value = 21.601941747572
$ws.Range("M16").Value = -31.228551818805
$ws.Range("N16").Value = -85.105908584169
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 89
$ws.Range("E17").Value = -26.966292134831
$ws.Range("F17").Value = 316
$ws.Range("G17").Value = 311
$ws.Range("H17").Value = 1.607717041800
$ws.Range("I17").Value = 3454
$ws.Range("J17").Value = 3403
$ws.Range("K17").Value = 1.498677637378
$ws.Range("L17").Value = 18.653383716935
$ws.Range("M17").Value = 26.520146520146
$ws.Range("N17").Value = -50.97232079489
$ws.Range("C18").Value = 43
$ws.Range("D18").Value = 56
$ws.Range("E18").Value = -23.214285714285
$ws.Range("F18").Value = 159
$ws.Range("G18").Value = 198
$ws.Range("H18").Value = -19.696969696969
$ws.Range("I18").Value = 1668
$ws.Range("J18").Value = 1955
$ws.Range("K18").Value = -14.680306905370
$ws.Range("L18").Value = 1.398176291793
$ws.Range("M18").Value = -35.697764070932
$ws.Range("N18").Value = -83.112281057001
$ws.Range("C19").Value = 120
$ws.Range("D19").Value = 139
$ws.Range("E19").Value = -13.669064748201
$ws.Range("F19").Value = 451
$ws.Range("G19").Value = 496
$ws.Range("H19").Value = -9.072580645161
$ws.Range("I19").Value = 4699
$ws.Range("J19").Value = 4843
$ws.Range("K19").Value = -2.973363617592
$ws.Range("L19").Value = 26.046137339055
$ws.Range("M19").Value = 36.757857974388
$ws.Range("N19").Value = -15.576715774344
$ws.Range("C20").Value = 42
$ws.Range("D20").Value = 30
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 162
$ws.Range("G20").Value = 154
$ws.Range("H20").Value = 5.194805194805
$ws.Range("I20").Value = 1498
$ws.Range("J20").Value = 1492
$ws.Range("K20").Value = 0.402144772117
$ws.Range("L20").Value = 21.197411003236
$ws.Range("M20").Value = 28.583690987124
$ws.Range("N20").Value = -80.364399003801
$ws.Range("C21").Value = 341
$ws.Range("D21").Value = 368
$ws.Range("E21").Value = -7.336956521739
$ws.Range("F21").Value = 1323
$ws.Range("G21").Value = 1378
$ws.Range("H21").Value = -3.991291727140
$ws.Range("I21").Value = 13564
$ws.Range("J21").Value = 14073
$ws.Range("K21").Value = -3.616854970510
$ws.Range("L21").Value = 18.680549479394
$ws.Range("M21").Value = 3.289674078586
$ws.Range("N21").Value = -69.490991700218
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -28.571428571428
$ws.Range("F22").Value = 13
$ws.Range("H22").Value = -45.833333333333
$ws.Range("I22").Value = 225
$ws.Range("J22").Value = 282
$ws.Range("K22").Value = -20.212765957446
$ws.Range("L22").Value = 9.756097560975
$ws.Range("M22").Value = -34.593023255814
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 28
$ws.Range("E23").Value = -14.285714285714
$ws.Range("F23").Value = 115
$ws.Range("G23").Value = 96
$ws.Range("H23").Value = 19.791666666666
$ws.Range("I23").Value = 1277
$ws.Range("J23").Value = 1234
$ws.Range("K23").Value = 3.484602917341
$ws.Range("L23").Value = 8.588435374149
$ws.Range("M23").Value = 33.159541188738
$ws.Range("C24").Value = 197
$ws.Range("D24").Value = 278
$ws.Range("E24").Value = -29.136690647482
$ws.Range("F24").Value = 812
$ws.Range("G24").Value = 1060
$ws.Range("H24").Value = -23.396226415094
$ws.Range("I24").Value = 10036
$ws.Range("J24").Value = 10889
$ws.Range("K24").Value = -7.833593534759
$ws.Range("L24").Value = 21.105345722215
$ws.Range("M24").Value = 18.335101992689
$ws.Range("C25").Value = 115
$ws.Range("D25").Value = 121
$ws.Range("E25").Value = -4.958677685950
$ws.Range("F25").Value = 453
$ws.Range("G25").Value = 394
$ws.Range("H25").Value = 14.974619289340
$ws.Range("I25").Value = 4969
$ws.Range("J25").Value = 4794
$ws.Range("K25").Value = 3.650396328744
$ws.Range("L25").Value = 30.351521511017
$ws.Range("M25").Value = -23.729854182655
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 30
$ws.Range("H26").Value = -11.764705882352
$ws.Range("I26").Value = 283
$ws.Range("J26").Value = 315
$ws.Range("K26").Value = -10.158730158730
$ws.Range("L26").Value = -8.116883116883
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 12
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 44
$ws.Range("G27").Value = 44
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 522
$ws.Range("J27").Value = 497
$ws.Range("K27").Value = 5.030181086519
$ws.Range("L27").Value = -6.451612903225
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 33.333333333333
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 30
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 195
$ws.Range("J28").Value = 292
$ws.Range("K28").Value = -33.219178082191
$ws.Range("L28").Value = -44.602272727272
$ws.Range("M28").Value = -55.275229357798
$ws.Range("N28").Value = -87.842892768079
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 33.333333333333
$ws.Range("F29").Value = 16
$ws.Range("G29").Value = 26
$ws.Range("H29").Value = -38.461538461538
$ws.Range("I29").Value = 165
$ws.Range("J29").Value = 245
$ws.Range("K29").Value = -32.653061224489
$ws.Range("L29").Value = -41.696113074204
$ws.Range("M29").Value = -53.521126760563
$ws.Range("N29").Value = -88.525730180806
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = -55.555555555555
$ws.Range("I30").Value = 51
$ws.Range("J30").Value = 70
$ws.Range("K30").Value = -27.142857142857
$ws.Range("L30").Value = 0

# C30 changes from a text/shared-string cell to a numeric cell; set number format to match column C
$ws.Range("C30").NumberFormat = "#,##0"
